# Add new daily data rows (44-49) to the worksheet, continuing the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$siteA = "四方坪站"
$siteB = "高岭站"

$rows = @(
    @{ Row = 44; Date = 45983; Site = $siteA; C = 9503.73;              D = 8599.1200000000008; E = 3149.8;   F = 405 },
    @{ Row = 45; Date = 45983; Site = $siteB; C = 5214.63;              D = 4505.25;             E = 1303.53; F = 188 },
    @{ Row = 46; Date = 45984; Site = $siteA; C = 8696.35;              D = 7592.61;             E = 2891.55; F = 369 },
    @{ Row = 47; Date = 45984; Site = $siteB; C = 3771.38;              D = 3262.13;             E = 946.75;  F = 150 },
    @{ Row = 48; Date = 45985; Site = $siteA; C = 8534.8700000000008;   D = 7645.7;              E = 2770.42; F = 370 },
    @{ Row = 49; Date = 45985; Site = $siteB; C = 3778.09;              D = 3116.36;             E = 1007.54; F = 143 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = $r.Date
    $ws.Cells.Item($i, 2).Value = $r.Site
    $ws.Cells.Item($i, 3).Value = $r.C
    $ws.Cells.Item($i, 4).Value = $r.D
    $ws.Cells.Item($i, 5).Value = $r.E
    $ws.Cells.Item($i, 6).Value = $r.F
}
